Write-Output "before"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$x = $s.Shapes.Item(999)
Write-Output "after"
Write-Output $x
